# Update benchmark: 2025-11-05 06:39:58 UTC
# Fills in previously-empty benchmark cells (and corrects one existing value)
# on the active worksheet of the workbook, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ŞANS OYUNLARI
$ws.Range("D2").Value = "23,81 TL - 23,81 TL"
$ws.Range("I2").Value = "18 TL - 18 TL"

# Row 3 - HESAPTAN EFT - Şube
$ws.Range("D3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4 - HESAPTAN EFT - ATM
$ws.Range("D4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5 - HESAPTAN EFT - Mobil
$ws.Range("D5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6 - DÜZENLİ EFT
$ws.Range("D6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("H6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("K6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 7 - KREDİ KARTINDAN FATURA ÖDEME
$ws.Range("D7").Value = "%1,6"

# Row 8 - HESAPTAN HAVALE - Şube
$ws.Range("D8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H8").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("K8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 9 - HESAPTAN HAVALE - ATM
$ws.Range("D9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H9").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("K9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 10 - HESAPTAN HAVALE - Mobil
$ws.Range("D10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H10").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("K10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 11 - DÜZENLİ HAVALE
$ws.Range("D11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("H11").Value = "3,05 TL - 6,1 TL - 76,18 TL"
$ws.Range("K11").Value = "3,05 TL - 6,09 TL - 76,17 TL"

# Row 12 - GİDEN SWIFT
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
$ws.Range("K12").Value = "WU: ,USD–; Diğer: 404,16 TL–3.403,42 TL"

# Row 13 - GELEN SWIFT
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"
$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"

# Row 14 - GİDEN SWIFT - Mobil
$ws.Range("D14").Value = "2.300 TL - 9.500 TL"
$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
$ws.Range("H14").Value = "2.100 TL - 4.300 TL"
$ws.Range("K14").Value = "914,14 TL - 4.265,98 TL"

# Row 24 - SENET TAHSİLE ALMA
$ws.Range("F24").Value = "457,14 TL"

# Row 25 - MUAMELESİZ SENET İADESİ
$ws.Range("F25").Value = "380,95 TL"
